# Apply updated TODC grade norms (raw score -> standard score) lookup tables
# through the IWR (grade-level) tabs, per commit "TODC grade norms through iwr".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("6.0-6.3")
$ws.Range("B2").Value = 57
$ws.Range("B3").Value = 72
$ws.Range("B5").Value = 94
$ws.Range("B6").Value = 103
$ws.Range("B7").Value = 112
$ws.Range("B8").Value = 120
$ws.Range("B9").Value = 127

$ws = $wb.Worksheets.Item("6.4-6.7")
$ws.Range("B2").Value = 55
$ws.Range("B3").Value = 69
$ws.Range("B4").Value = 81
$ws.Range("B5").Value = 92
$ws.Range("B6").Value = 101
$ws.Range("B7").Value = 109
$ws.Range("B8").Value = 117
$ws.Range("B9").Value = 124
$ws.Range("B10").Value = 130

$ws = $wb.Worksheets.Item("6.8-6.11")
$ws.Range("B2").Value = 53
$ws.Range("B3").Value = 67
$ws.Range("B4").Value = 79
$ws.Range("B6").Value = 99
$ws.Range("B7").Value = 107
$ws.Range("B8").Value = 115
$ws.Range("B9").Value = 122
$ws.Range("B10").Value = 129

$ws = $wb.Worksheets.Item("7.0-7.3")
$ws.Range("B2").Value = 51
$ws.Range("B3").Value = 66
$ws.Range("B4").Value = 77
$ws.Range("B5").Value = 87
$ws.Range("B6").Value = 97
$ws.Range("B7").Value = 105
$ws.Range("B8").Value = 113
$ws.Range("B9").Value = 120
$ws.Range("B10").Value = 126

$ws = $wb.Worksheets.Item("7.4-7.7")
$ws.Range("B2").Value = 49
$ws.Range("B3").Value = 64
$ws.Range("B4").Value = 76
$ws.Range("B9").Value = 118

$ws = $wb.Worksheets.Item("7.8-7.11")
$ws.Range("B2").Value = 47
$ws.Range("B3").Value = 62
$ws.Range("B4").Value = 74
$ws.Range("B5").Value = 84

$ws = $wb.Worksheets.Item("8.0-8.5")
$ws.Range("B2").Value = 45
$ws.Range("B3").Value = 60
$ws.Range("B4").Value = 72
$ws.Range("B5").Value = 82
$ws.Range("B6").Value = 91
$ws.Range("B7").Value = 99
$ws.Range("B8").Value = 106
$ws.Range("B9").Value = 113
$ws.Range("B10").Value = 120
$ws.Range("B11").Value = 126

$ws = $wb.Worksheets.Item("8.6-8.11")
$ws.Range("B2").Value = 43
$ws.Range("B3").Value = 58
$ws.Range("B4").Value = 70
$ws.Range("B5").Value = 80
$ws.Range("B6").Value = 89
$ws.Range("B7").Value = 97
$ws.Range("B8").Value = 104
$ws.Range("B9").Value = 111
$ws.Range("B10").Value = 117
$ws.Range("B11").Value = 124
$ws.Range("B12").Value = 129

$ws = $wb.Worksheets.Item("9.0-9.5")
$ws.Range("B2").Value = 40
$ws.Range("B3").Value = 56
$ws.Range("B4").Value = 68
$ws.Range("B5").Value = 78
$ws.Range("B6").Value = 87
$ws.Range("B7").Value = 95
$ws.Range("B8").Value = 102
$ws.Range("B9").Value = 109
$ws.Range("B10").Value = 115
$ws.Range("B11").Value = 121
$ws.Range("B12").Value = 127

$ws = $wb.Worksheets.Item("9.6-9.11")
$ws.Range("B2").Value = 40
$ws.Range("B3").Value = 54
$ws.Range("B4").Value = 66
$ws.Range("B5").Value = 76
$ws.Range("B6").Value = 85
$ws.Range("B7").Value = 93
$ws.Range("B8").Value = 100
$ws.Range("B9").Value = 107
$ws.Range("B10").Value = 113
$ws.Range("B11").Value = 119
$ws.Range("B12").Value = 125

$ws = $wb.Worksheets.Item("10.0-10.5")
$ws.Range("B2").Value = 40
$ws.Range("B3").Value = 52
$ws.Range("B4").Value = 64
$ws.Range("B5").Value = 74
$ws.Range("B6").Value = 83
$ws.Range("B7").Value = 91
$ws.Range("B8").Value = 99
$ws.Range("B9").Value = 105
$ws.Range("B10").Value = 112
$ws.Range("B11").Value = 117
$ws.Range("B12").Value = 123
$ws.Range("B13").Value = 128

$ws = $wb.Worksheets.Item("10.6-10.11")
$ws.Range("B2").Value = 40
$ws.Range("B3").Value = 50
$ws.Range("B4").Value = 63
$ws.Range("B5").Value = 73
$ws.Range("B6").Value = 82
$ws.Range("B7").Value = 90
$ws.Range("B8").Value = 97
$ws.Range("B9").Value = 104
$ws.Range("B10").Value = 110
$ws.Range("B11").Value = 116
$ws.Range("B12").Value = 121
$ws.Range("B13").Value = 127

$ws = $wb.Worksheets.Item("11.0-11.5")
$ws.Range("B2").Value = 40
$ws.Range("B3").Value = 48
$ws.Range("B4").Value = 61
$ws.Range("B5").Value = 72
$ws.Range("B6").Value = 80
$ws.Range("B7").Value = 88
$ws.Range("B8").Value = 96
$ws.Range("B9").Value = 102
$ws.Range("B10").Value = 109
$ws.Range("B11").Value = 114
$ws.Range("B12").Value = 120
$ws.Range("B13").Value = 125

$ws = $wb.Worksheets.Item("11.6-11.11")
$ws.Range("B2").Value = 40
$ws.Range("B3").Value = 47
$ws.Range("B4").Value = 60
$ws.Range("B5").Value = 70
$ws.Range("B6").Value = 79
$ws.Range("B7").Value = 87
$ws.Range("B8").Value = 95
$ws.Range("B9").Value = 101
$ws.Range("B10").Value = 107
$ws.Range("B11").Value = 113
$ws.Range("B12").Value = 119
$ws.Range("B13").Value = 124
$ws.Range("B14").Value = 129

$ws = $wb.Worksheets.Item("12.0-12.5")
$ws.Range("B3").Value = 45
$ws.Range("B4").Value = 58
$ws.Range("B5").Value = 69
$ws.Range("B6").Value = 78
$ws.Range("B7").Value = 86
$ws.Range("B8").Value = 93
$ws.Range("B9").Value = 100
$ws.Range("B10").Value = 106
$ws.Range("B11").Value = 112
$ws.Range("B12").Value = 117
$ws.Range("B13").Value = 122
$ws.Range("B14").Value = 127

$ws = $wb.Worksheets.Item("12.6-12.11")
$ws.Range("B3").Value = 43
$ws.Range("B4").Value = 57
$ws.Range("B5").Value = 68
$ws.Range("B6").Value = 77
$ws.Range("B7").Value = 85
$ws.Range("B8").Value = 93
$ws.Range("B9").Value = 99
$ws.Range("B10").Value = 105
$ws.Range("B11").Value = 111
$ws.Range("B12").Value = 116
$ws.Range("B13").Value = 121
$ws.Range("B14").Value = 126

$ws = $wb.Worksheets.Item("13.0-13.11")
$ws.Range("B3").Value = 41
$ws.Range("B4").Value = 55
$ws.Range("B5").Value = 66
$ws.Range("B6").Value = 76
$ws.Range("B7").Value = 84
$ws.Range("B8").Value = 91
$ws.Range("B9").Value = 98
$ws.Range("B10").Value = 104
$ws.Range("B11").Value = 110
$ws.Range("B12").Value = 115
$ws.Range("B13").Value = 120
$ws.Range("B14").Value = 125
$ws.Range("B15").Value = 129

$ws = $wb.Worksheets.Item("14.0-14.11")
$ws.Range("B3").Value = 40
$ws.Range("B4").Value = 53
$ws.Range("B5").Value = 65
$ws.Range("B6").Value = 74
$ws.Range("B8").Value = 90
$ws.Range("B10").Value = 103
$ws.Range("B11").Value = 108
$ws.Range("B12").Value = 113
$ws.Range("B13").Value = 118
$ws.Range("B14").Value = 123
$ws.Range("B15").Value = 127

$ws = $wb.Worksheets.Item("15.0-16.11")
$ws.Range("B3").Value = 40
$ws.Range("B4").Value = 50
$ws.Range("B5").Value = 63
$ws.Range("B6").Value = 73
$ws.Range("B7").Value = 81
$ws.Range("B9").Value = 95
$ws.Range("B10").Value = 101
$ws.Range("B15").Value = 125

$ws = $wb.Worksheets.Item("17.0-18.11")
$ws.Range("B3").Value = 40
$ws.Range("B4").Value = 46
$ws.Range("B5").Value = 61
$ws.Range("B6").Value = 72
$ws.Range("B9").Value = 95
$ws.Range("B11").Value = 106
$ws.Range("B16").Value = 126
